$wb = $excel.ActiveWorkbook

# Fix sheet name typos
$wb.Worksheets.Item("Bevereges").Name = "Beverages"
$wb.Worksheets.Item("Hygeine & Personal Care").Name = "Hygene & Personal Care"

# Move the active/selected tab from " Baby & Childcare" (last sheet) to
# "Household Essentials" (the previous sheet), matching the updated
# activeTab / tabSelected state in the workbook.
$wb.Worksheets.Item("Household Essentials").Activate()
